# Trade #14 closed at 2026-02-17 07:53:58 - unknown UNKNOWN +0.000%
#
# Updates the Summary + Strategy Status roll-up numbers to reflect the new
# closed trade, and appends the new trade row to both the "All Trades" and
# "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.94   # Current Capital
$summary.Range("B4").Value = -0.06    # Total P&L $
$summary.Range("B5").Value = -0.09    # Total P&L %
$summary.Range("B6").Value = 14        # Total Trades
$summary.Range("B8").Value = 8         # Losing Trades
$summary.Range("B9").Value = 35.71    # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.94     # Capital
$status.Range("D4").Value = 14         # Trades
$status.Range("E4").Value = -0.06     # P&L $
$status.Range("F4").Value = -0.06     # P&L %
$status.Range("G4").Value = 35.71     # Win Rate %

# ---------------------------------------------------------------------
# Helper that appends the new trade-#14 row to a trade-log sheet.
# The Date column ("2026-02-17") looks like a date, so the cell has to be
# pre-formatted as Text before the value is written - otherwise Excel's
# auto-detection would silently convert it into a date serial number.
# ---------------------------------------------------------------------
function Add-TradeRow14($ws) {
    $ws.Range("A15").Value = 14
    $ws.Range("B15").NumberFormat = "@"
    $ws.Range("B15").Value = "2026-02-17"
    $ws.Range("C15").Value = "07:53:51"
    $ws.Range("D15").Value = "MarketMaking"
    $ws.Range("E15").Value = "DOWN"
    $ws.Range("F15").Value = 0.3
    $ws.Range("G15").Value = 0.28
    $ws.Range("H15").Value = "CLOSED"
    $ws.Range("I15").Value = -6.6667
    $ws.Range("J15").Value = -0.02
    $ws.Range("K15").Value = 99.94
    $ws.Range("L15").Value = 0
    $ws.Range("M15").Value = 0
    $ws.Range("N15").Value = 0.6
    $ws.Range("O15").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P15").Value = "early_exit"
    $ws.Range("Q15").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow14 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow14 $marketMaking
